$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.611.36'
$ws.Range('E2').Value = '  +2.56%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.993.99'
$ws.Range('E3').Value = '  +6.14%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '328.10'
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4685'
$ws.Range('E7').Value = '  +1.49%  '
$ws.Range('E8').Value = '  +1.88%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.70'
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08035'
$ws.Range('E10').Value = '  +2.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.002'
$ws.Range('E11').Value = '  +1.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.94'
$ws.Range('E12').Value = '  +5.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.997.34'
$ws.Range('E13').Value = '  +8.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.247'
$ws.Range('E14').Value = '  +3.60%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.875'
$ws.Range('E15').Value = '  +3.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.07124'
$ws.Range('E16').Value = '  +1.99%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '89.03'
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.003'
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001004'
$ws.Range('E19').Value = '  +0.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.44'
$ws.Range('E20').Value = '  +2.81%  '
$ws.Range('E21').Value = '  -0.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '29.627.76'
$ws.Range('E22').Value = '  +2.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.564'
$ws.Range('E23').Value = '  +5.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.26'
$ws.Range('E24').Value = '  +2.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.104'
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.09'
$ws.Range('E26').Value = '  +1.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.71'
$ws.Range('E27').Value = '  +1.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.982'
$ws.Range('E28').Value = '  +1.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '120.32'
$ws.Range('E29').Value = '  +2.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.952'
$ws.Range('E30').Value = '  +2.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09462'
$ws.Range('E31').Value = '  +1.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9184'
$ws.Range('E32').Value = '  +1.94%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.362'
$ws.Range('E33').Value = '  +3.20%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.289'
$ws.Range('E34').Value = '  +0.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.214'
$ws.Range('E35').Value = '  -1.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05863'
$ws.Range('E36').Value = '  +2.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.175'
$ws.Range('E37').Value = '  +0.58%  '
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.000003384'
$ws.Range('E38').Value = '  +70.46%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02128'
$ws.Range('E39').Value = '  +2.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.924'
$ws.Range('E40').Value = '  +4.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5798'
$ws.Range('E41').Value = '  +2.50%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1828'
$ws.Range('E42').Value = '  +3.06%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '9.908'
$ws.Range('E43').Value = '  +2.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.809'
$ws.Range('E44').Value = '  +10.98%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.04'
$ws.Range('E45').Value = '  +0.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5411'
$ws.Range('E46').Value = '  +1.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.212'
$ws.Range('E47').Value = '  -1.55%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.874'
$ws.Range('E48').Value = '  +1.70%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06972'
$ws.Range('E49').Value = '  -0.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '114.39'
$ws.Range('E50').Value = '  +1.71%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3085'
$ws.Range('E51').Value = '  +8.14%  '
